# Update the "Förändrad" (Changed) date column (C) for rows 2-20
# from 45243 (2023-11-13) to 45244 (2023-11-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$ws.Range("C2:C20").Value = 45244
